$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.696681380271912
$ws.Range("B1").Value = 2.546590805053711
$ws.Range("C1").Value = 3.157073736190796
$ws.Range("D1").Value = 2.389169692993164
$ws.Range("E1").Value = 0.460883766412735
